$d = $word.ActiveDocument

# Helper: replace the text of a whole paragraph (by 1-based index) while
# keeping the paragraph mark (and therefore the run/paragraph formatting
# of the trailing run) intact.
function Set-ParagraphText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1)
    $r.Text = $newText
}

# 1) Title heading, appears twice with identical old/new text -> use a
#    global Find/Replace so both the Heading1 at the top and the bold
#    "Play ..." run near the bottom get updated together.
$d.Content.Find.Execute(
    "Play Cat Wilde and the Doom of Dead Free | Review", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Play Cat Wilde and the Doom of Dead for Free", 2)

# 2) "What we like" bullet list updates
Set-ParagraphText 43 "Free spins feature with special expanding symbol"
Set-ParagraphText 44 "Payouts of up to 5,000x stake per spin"
Set-ParagraphText 45 "Excellent design and high-quality graphics"

# 3) "What we don't like" bullet list updates
Set-ParagraphText 47 "High volatility may not be suited for low variance slot players"
Set-ParagraphText 48 "Limited wagering options"

# 4) Meta description (italic) paragraph near the very end
Set-ParagraphText 50 "Read our review of Cat Wilde and the Doom of Dead, a highly volatile slot game with expanding wilds and free spins. Play for free now."
